# Eurovision song contest matrix: the sheet is a "from country" x "to
# country" points matrix (rows 4-47 = "from" countries in column A, columns
# B:AQ = "to" countries, header in row 2). The diagonal cells - where a
# country would have voted for itself - were incorrectly stored as a
# literal 0. This fixes the diagonal by blanking those self-vote cells out
# instead of leaving a numeric 0 in them.
#
# Note: two "from" countries (Bosnia & Herzegovina, row 11; Montenegro, row
# 32) have no matching "to" column, so they have no diagonal cell to touch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$diagonalCells = @(
    "B4", "C5", "D6", "E7", "F8", "G9", "H10",
    "I12", "J13", "K14", "L15", "M16", "N17", "O18", "P19", "Q20", "R21",
    "S22", "T23", "U24", "V25", "W26", "X27", "Y28", "Z29", "AA30", "AB31",
    "AC33", "AD34", "AE35", "AF36", "AG37", "AH38", "AI39", "AJ40", "AK41",
    "AL42", "AM43", "AN44", "AO45", "AP46", "AQ47"
)

foreach ($cellRef in $diagonalCells) {
    $ws.Range($cellRef).ClearContents()
}
